$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 330.33334
$ws.Range("I33").Value = 306.66666
$ws.Range("J33").Value = 354
$ws.Range("K33").Value = 306.66666
$ws.Range("L33").Value = 354
$ws.Range("M33").Value = -77.66665999999998
$ws.Range("N33").Value = -812
$ws.Range("H44").Value = 4982.706
$ws.Range("I44").Value = 2246
$ws.Range("J44").Value = 10000
$ws.Range("K44").Value = 2246
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = -1784
$ws.Range("N44").Value = -10924
$ws.Range("H64").Value = 3619.1904
$ws.Range("I64").Value = 3160
$ws.Range("J64").Value = 4036.6365
$ws.Range("K64").Value = 3160
$ws.Range("L64").Value = 4036.6365
$ws.Range("M64").Value = -2912
$ws.Range("N64").Value = -4532.636500000001
$ws.Range("H67").Value = 3619.1904
$ws.Range("I67").Value = 3160
$ws.Range("J67").Value = 4036.6365
$ws.Range("K67").Value = 3160
$ws.Range("L67").Value = 4036.6365
$ws.Range("M67").Value = -2302
$ws.Range("N67").Value = -5752.636500000001
$ws.Range("H74").Value = 7545.5713
$ws.Range("I74").Value = 3401.5
$ws.Range("J74").Value = 9203.200000000001
$ws.Range("K74").Value = 3401.5
$ws.Range("L74").Value = 9203.200000000001
$ws.Range("M74").Value = -2465.5
$ws.Range("N74").Value = -11075.2
$ws.Range("H77").Value = 7545.5713
$ws.Range("I77").Value = 3401.5
$ws.Range("J77").Value = 9203.200000000001
$ws.Range("K77").Value = 17007.5
$ws.Range("L77").Value = 46016
$ws.Range("M77").Value = -12327.5
$ws.Range("N77").Value = -55376

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8159.1113
$ws.Range("J3").Value = 7833.3335
$ws.Range("L3").Value = 23500.0005
$ws.Range("N3").Value = -23724.0005
$ws.Range("H59").Value = 838
$ws.Range("I59").Value = 297.5
$ws.Range("K59").Value = 892.5
$ws.Range("M59").Value = -352.5
$ws.Range("H70").Value = 933.3333
$ws.Range("I70").Value = 400
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 1200
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -885
$ws.Range("N70").Value = -6630
$ws.Range("H73").Value = 933.3333
$ws.Range("I73").Value = 400
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 1200
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -108
$ws.Range("N73").Value = -8184
$ws.Range("H94").Value = 1833.3334
$ws.Range("I94").Value = 1750
$ws.Range("K94").Value = 5250
$ws.Range("M94").Value = -4574
$ws.Range("H96").Value = 8740
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 8740
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 26220
$ws.Range("M96").Value = ""
$ws.Range("N96").Value = -30338
$ws.Range("H97").Value = 788.8
$ws.Range("I97").Value = 663.3333
$ws.Range("J97").Value = 977
$ws.Range("K97").Value = 1989.9999
$ws.Range("L97").Value = 2931
$ws.Range("M97").Value = -1493.9999
$ws.Range("N97").Value = -3923
$ws.Range("H103").Value = 730.8
$ws.Range("I103").Value = 32
$ws.Range("J103").Value = 1196.6666
$ws.Range("K103").Value = 96
$ws.Range("L103").Value = 3589.9998
$ws.Range("M103").Value = 783
$ws.Range("N103").Value = -5347.9998
$ws.Range("H104").Value = 44621.668
$ws.Range("J104").Value = 66432.5
$ws.Range("L104").Value = 199297.5
$ws.Range("N104").Value = -204539.5
$ws.Range("H105").Value = 4602.231
$ws.Range("I105").Value = 2644.4443
$ws.Range("J105").Value = 9007.25
$ws.Range("K105").Value = 7933.3329
$ws.Range("L105").Value = 27021.75
$ws.Range("M105").Value = -5312.3329
$ws.Range("N105").Value = -32263.75
$ws.Range("H109").Value = 2653.2693
$ws.Range("I109").Value = 691.6923
$ws.Range("J109").Value = 4614.846
$ws.Range("K109").Value = 2075.0769
$ws.Range("L109").Value = 13844.538
$ws.Range("M109").Value = -1035.0769
$ws.Range("N109").Value = -15924.538
$ws.Range("H115").Value = 5686.087
$ws.Range("I115").Value = 800
$ws.Range("J115").Value = 5908.1816
$ws.Range("K115").Value = 2400
$ws.Range("L115").Value = 17724.5448
$ws.Range("M115").Value = -1225
$ws.Range("N115").Value = -20074.5448
$ws.Range("H118").Value = 22608
$ws.Range("I118").Value = 482.25
$ws.Range("J118").Value = 111111
$ws.Range("K118").Value = 1446.75
$ws.Range("L118").Value = 333333
$ws.Range("M118").Value = -203.75
$ws.Range("N118").Value = -335819
$ws.Range("H121").Value = 7937165
$ws.Range("I121").Value = 293.5
$ws.Range("J121").Value = 27779344
$ws.Range("K121").Value = 880.5
$ws.Range("L121").Value = 83338032
$ws.Range("M121").Value = 429.5
$ws.Range("N121").Value = -83340652
$ws.Range("H122").Value = 18522828
$ws.Range("I122").Value = 38461812
$ws.Range("J122").Value = 8059.0713
$ws.Range("K122").Value = 346156308
$ws.Range("L122").Value = 72531.64169999999
$ws.Range("M122").Value = -346153858
$ws.Range("N122").Value = -77431.64169999999
$ws.Range("H123").Value = 300
$ws.Range("I123").Value = 300
$ws.Range("K123").Value = 900
$ws.Range("M123").Value = 1550
$ws.Range("H126").Value = 2836.842
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 2985.7144
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 8957.143199999999
$ws.Range("M126").Value = 1640
$ws.Range("N126").Value = -18837.1432

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3174.4407
$ws.Range("I132").Value = 1274.4419
$ws.Range("J132").Value = 8280.6875
$ws.Range("K132").Value = 3823.3257
$ws.Range("L132").Value = 24842.0625
$ws.Range("M132").Value = -1293.3257
$ws.Range("N132").Value = -29902.0625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""
